$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated predictions for column B (row 2..8) -- rerun "con todas las variables sin depuracion"
$ws.Cells.Item(2,2).Value = 34988.91796875
$ws.Cells.Item(3,2).Value = 34983.72265625
$ws.Cells.Item(4,2).Value = 34958.0625
$ws.Cells.Item(5,2).Value = 34710.375
$ws.Cells.Item(6,2).Value = 34692.16015625
$ws.Cells.Item(7,2).Value = 34723.28125
$ws.Cells.Item(8,2).Value = 36375.9609375
